$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Group F3E qsQDUv"
$ws.Range("A3").Value = "Group Lf7G 8 c  "
$ws.Range("A4").Value = "Group b fF TO Va"
$ws.Range("A5").Value = "Group FRS zalFlX"
$ws.Range("A6").Value = "Group 8fUEBFhb0P"
